$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds Property/Value pairs.
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 = "Name" property -> set its value (was empty)
$ws.Range("B4").Value = "AnneeuniversitaireVs"

# Row 8 = "Date" property -> update the generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
